$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("StrategyDictionaries")

# --- Sheet1: new HA_VWAP test rows ---

# Row 2 (Test #1): change optional strategy settings to the corrected
# (well-formed) JSON with DistVWAP_PCT = 0.0
$ws1.Range("K2").Value = '{"EMA": 200, "DistVWAP_PCT": 0.0, "NB_SIGNALS": 2}'

# Row 3 (Test #2): pair switches from ETHUSDT back to BTCUSDT ...
$ws1.Range("C3").Value = "BTCUSDT"
# ... and gets the HA_VWAP settings with DistVWAP_PCT = 0.05
$ws1.Range("K3").Value = '{"EMA": 200, "DistVWAP_PCT": 0.05, "NB_SIGNALS": 2}'

# --- StrategyDictionaries: fix the HA_VWAP reference dictionary text ---
$ws3.Range("B5").Value = '{"EMA": 200, "DistVWAP_PCT": 0.05, "NB_SIGNALS": 2}'

# --- View state: Sheet1 was being worked on (scrolled right, K3 selected) ---
$ws1.Activate()
$ws1.Range("K3").Select()

# --- Final active sheet ends up being StrategyDictionaries ---
$ws3.Activate()
$ws3.Range("B11").Select()
